$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 7; this shifts existing rows 7-23 down to 8-24,
# preserving all their data/formatting (matches the diff's row re-numbering).
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with a fresh weekly price record.
$ws.Cells.Item(7,1).Value = 1
$ws.Cells.Item(7,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(7,3).Value = "Arica y Parinacota"
$ws.Cells.Item(7,4).Value = 44525
$ws.Cells.Item(7,4).NumberFormat = $ws.Cells.Item(8,4).NumberFormat
$ws.Cells.Item(7,5).Value = 15
$ws.Cells.Item(7,6).Value = 100112044
$ws.Cells.Item(7,7).Value = "Perejil"
$ws.Cells.Item(7,8).Value = "Sin especificar"
$ws.Cells.Item(7,9).Value = "Primera"
$ws.Cells.Item(7,10).Value = 300
$ws.Cells.Item(7,11).Value = 1400
$ws.Cells.Item(7,12).Value = 1500
$ws.Cells.Item(7,13).Value = 1450
$ws.Cells.Item(7,14).Value = "$/atado 1,5 a 2 kilos"
$ws.Cells.Item(7,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(7,16).Value = 725
$ws.Cells.Item(7,17).Value = 2
$ws.Cells.Item(7,18).Value = "Hortaliza"
